$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the "#" numbering column (A2:A12) for the existing findings.
# ---------------------------------------------------------------------------
for ($i = 2; $i -le 12; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Row 13 previously had no value in column A - give it number 12 and make
# sure it picks up the same bordered style used by the rows above it.
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 12

# Row 13's Impact/Description/Solution/Level cells used the "no border"
# style; align them with the bordered style used by row 12 just above.
$ws.Range("C12:F12").Copy()
$ws.Range("C13:F13").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Append the new findings as rows 14-18.
# ---------------------------------------------------------------------------

# Row 14
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 13

$ws.Range("B12").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = "Quitar Impresión de matriz (DDC)"

$ws.Range("C12:F12").Copy()
$ws.Range("C14:F14").PasteSpecial(-4122)
$ws.Range("C14").Value = "Medio"
$ws.Range("D14").Value = "Se permite imprimir desde  la vista de DDC"
$ws.Range("E14").Value = "Un usario que solo se le permite vista de DDC Puede realizar una impresión de DDC"
$ws.Range("F14").Value = "Error de codigo fuente"

$ws.Range("G12").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = "Abierta"

$ws.Rows.Item(14).RowHeight = 31.5

# Row 15
$ws.Range("A2").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 14

$ws.Range("B12").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "Revision de DDC en apartado de impresión"

$ws.Range("C12:F12").Copy()
$ws.Range("C15:F15").PasteSpecial(-4122)
$ws.Range("C15").Value = "Medio"
$ws.Range("D15").Value = "Error de formato DDC"
$ws.Range("E15").Value = " La impresión no se generá en una sola pagina, la impresión no contiene una vista estetica para impresión"
$ws.Range("F15").Value = "Error de codigo fuente "

$ws.Range("G12").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = "Abierta"

$ws.Rows.Item(15).RowHeight = 47.25

# Row 16
$ws.Range("A2").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 15

$ws.Range("B12").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Revision de los CBX en blanco"

$ws.Range("C12:F12").Copy()
$ws.Range("C16:F16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Alta"
$ws.Range("D16").Value = "Atributo no seleccionado"
$ws.Range("E16").Value = "El Atributo que se selecciona en el Combobox no aparece como seleccionado a nivel de interfaz, si logra tomar el valor pero no es mostrado en la UI"
$ws.Range("F16").Value = "Error de codigo fuente"

$ws.Range("G12").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = "Abierta"

$ws.Rows.Item(16).RowHeight = 47.25

# Row 17
$ws.Range("A2").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 16

$ws.Range("B12").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "Administación para Calibración de Factores matriz"

$ws.Range("C12:F12").Copy()
$ws.Range("C17:F17").PasteSpecial(-4122)
$ws.Range("C17").Value = "Medio"
$ws.Range("D17").Value = "Error de rotulacion"
$ws.Range("E17").Value = "El apartado no ha sido rotulado de forma correcta y se desconoce  el subconjuto de datos a los que hace referencia"
$ws.Range("F17").Value = "Error de codigo fuente -Capa presentacion"

$ws.Range("G12").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G17").Value = "Abierta"

$ws.Rows.Item(17).RowHeight = 47.25

# Row 18 - no "#" value this time; B18 uses a new left-aligned bordered style.
$ws.Range("D12").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").HorizontalAlignment = -4131
$ws.Range("B18").Value = "Cambiar listaInternas a wv_listas consol"

$ws.Range("C12:F12").Copy()
$ws.Range("C18:F18").PasteSpecial(-4122)
$ws.Range("C18").Value = "Medio"
$ws.Range("D18").Value = "Agregar base de datos de FBI, ONU, Listas Internas; (pep_historial)"
$ws.Range("E18").Value = "Se deberá de crear una tabla consolidad que contenga las base de datos FBI, PEP (ACTIVOS), ENGEL, INTERPOL, SANCION BANCO MUDIAL"
$ws.Range("F18").Value = "Integracion incompleta"

$ws.Range("G12").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("G18").Value = "Abierta"

$ws.Rows.Item(18).RowHeight = 47.25

# ---------------------------------------------------------------------------
# 3) Update the view: scroll back to the top and select the next empty row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("B20").Select()
